# Rajasthan Royals / Ben Stokes sheet: correct the per-innings
# runs/balls/fours/sixes figures in rows 2-9 (columns C:F) so that each
# row holds the right stat line. Cells are stored as text (the sheet
# already flags "number stored as text"), so force NumberFormat "@"
# before writing each value to keep them text instead of letting Excel
# auto-coerce them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new C (runs), D (balls), E (fours), F (sixes) values
$updates = @(
    @{ Row = 2; C = "15";  D = "19"; E = "2";  F = "0" },
    @{ Row = 3; C = "50";  D = "26"; E = "6";  F = "3" },
    @{ Row = 4; C = "30";  D = "32"; E = "2";  F = "0" },
    @{ Row = 5; C = "107"; D = "60"; E = "14"; F = "3" },
    @{ Row = 6; C = "41";  D = "35"; E = "6";  F = "0" },
    @{ Row = 7; C = "5";   D = "6";  E = "1";  F = "0" },
    @{ Row = 8; C = "18";  D = "11"; E = "2";  F = "1" },
    @{ Row = 9; C = "19";  D = "11"; E = "3";  F = "0" }
)

foreach ($update in $updates) {
    foreach ($col in "C", "D", "E", "F") {
        $cellRef = "$col$($update.Row)"
        $newValue = $update[$col]
        $cell = $ws.Range($cellRef)

        # Only touch cells whose value actually changes, and keep them
        # typed as text (matches the original "stored as text" cells).
        if ($cell.Text -ne $newValue) {
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
        }
    }
}
